# Weekly data refresh: insert one new price observation as a new row
# right above the existing row 455, shifting all subsequent rows down
# by one (dimension grows from A1:R575 to A1:R576).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 455; Excel shifts rows 455:575 down to 456:576
# and copies formatting (incl. the date number format on column D) from
# the row above, matching the existing style index used throughout col D.
$ws.Rows("455:455").Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A455").Value = 10
$ws.Range("B455").Value = "Vega Modelo de Temuco"
$ws.Range("C455").Value = "La Araucanía"
$ws.Range("D455").Value = 44932
$ws.Range("E455").Value = 9
$ws.Range("F455").Value = 100112023
$ws.Range("G455").Value = "Brócoli"
$ws.Range("H455").Value = "Sin especificar"
$ws.Range("I455").Value = "Primera"
$ws.Range("J455").Value = 200
$ws.Range("K455").Value = 1200
$ws.Range("L455").Value = 1200
$ws.Range("M455").Value = 1200
$ws.Range("N455").Value = "$/unidad"
$ws.Range("O455").Value = "Provincia de Cautín"
$ws.Range("P455").Value = 1200
$ws.Range("Q455").Value = 1
$ws.Range("R455").Value = "Hortaliza"
